$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as plain text, matching the
# original inline-string storage (avoids numeric auto-conversion / trailing-zero loss).

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '68.459.66'
$ws.Range('E2').Value = '  -2.73%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.438.31'
$ws.Range('E3').Value = '  -5.00%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '571.57'
$ws.Range('E5').Value = '  -5.02%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '188.91'
$ws.Range('E6').Value = '  -4.03%  '
$ws.Range('E7').Value = '  -3.72%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.426.30'
$ws.Range('E8').Value = '  -5.01%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('E10').Value = '  -5.82%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.611'
$ws.Range('E11').Value = '  -5.46%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '50.70'
$ws.Range('E12').Value = '  -4.80%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000282'
$ws.Range('E13').Value = '  -7.56%  '
$ws.Range('E14').Value = '  -5.74%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.979.73'
$ws.Range('E15').Value = '  -5.13%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '630.43'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '68.287.69'
$ws.Range('E17').Value = '  -3.07%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.446.45'
$ws.Range('E18').Value = '  -4.57%  '
$ws.Range('E19').Value = '  -2.56%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.17'
$ws.Range('E20').Value = '  -5.76%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.97'
$ws.Range('E21').Value = '  -5.57%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.933'
$ws.Range('E22').Value = '  -6.55%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '17.67'
$ws.Range('E23').Value = '  -2.57%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.30'
$ws.Range('E24').Value = '  +1.83%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '98.34'
$ws.Range('E25').Value = '  -4.56%  '
$ws.Range('E26').Value = '  -8.30%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.06'
$ws.Range('E27').Value = '  +1.89%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.81'
$ws.Range('E28').Value = '  -5.91%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.74'
$ws.Range('E29').Value = '  -8.28%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '9.10'
$ws.Range('E31').Value = '  -5.13%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.11'
$ws.Range('E32').Value = '  -12.22%  '
$ws.Range('E33').Value = '  -8.88%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '11.48'
$ws.Range('E34').Value = '  -6.39%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '60.50'
$ws.Range('E35').Value = '  -4.45%  '
$ws.Range('E36').Value = '  -7.70%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.06%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.623.52'
$ws.Range('E38').Value = '  -7.86%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0₃0777'
$ws.Range('E39').Value = '  -12.62%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '500.17'
$ws.Range('E40').Value = '  -4.13%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.86'
$ws.Range('E41').Value = '  -6.51%  '
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('E43').Value = '  -6.31%  '
$ws.Range('E44').Value = '  -3.45%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '33.91'
$ws.Range('E45').Value = '  -7.57%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.45'
$ws.Range('E46').Value = '  +66.83%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0433'
$ws.Range('E47').Value = '  -6.18%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.32'
$ws.Range('E48').Value = '  -5.45%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.77'
$ws.Range('E49').Value = '  -4.57%  '
$ws.Range('E50').Value = '  -5.08%  '
$ws.Range('E51').Value = '  -0.36%  '
